$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.028965631712938
$ws.Range("D2").Value = 1.037227525576902
$ws.Range("E2").Value = 0.9926147277508489
$ws.Range("F2").Value = 1.044210805610972
$ws.Range("I2").Value = 1.032231077692212
$ws.Range("J2").Value = 1.034115024638117
$ws.Range("K2").Value = 1.040019011186655
$ws.Range("L2").Value = 0.9955398523336033
$ws.Range("M2").Value = 1.04698249520968
$ws.Range("N2").Value = 1.035583586587775

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.029940688537354
$ws.Range("D3").Value = 1.037987402974559
$ws.Range("E3").Value = 0.9936372048519304
$ws.Range("F3").Value = 1.04516135259206
$ws.Range("I3").Value = 1.032379384310972
$ws.Range("J3").Value = 1.034730782812199
$ws.Range("K3").Value = 1.040588785002228
$ws.Range("L3").Value = 0.9963617723202692
$ws.Range("M3").Value = 1.047743870258755
$ws.Range("N3").Value = 1.036200219209072

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.030571906238008
$ws.Range("D4").Value = 1.038479156597969
$ws.Range("E4").Value = 0.9942998659930995
$ws.Range("F4").Value = 1.045776967652096
$ws.Range("I4").Value = 1.032473993568972
$ws.Range("J4").Value = 1.035128920655808
$ws.Range("K4").Value = 1.04095686788792
$ws.Range("L4").Value = 0.9968940712668345
$ws.Range("M4").Value = 1.048236456269194
$ws.Range("N4").Value = 1.036598922454086

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.030837339138631
$ws.Range("D5").Value = 1.038685903502808
$ws.Range("E5").Value = 0.9945786998346017
$ws.Range("F5").Value = 1.04603590218602
$ws.Range("I5").Value = 1.03251344243277
$ws.Range("J5").Value = 1.03529622549898
$ws.Range("K5").Value = 1.041111465669839
$ws.Range("L5").Value = 0.997117960005301
$ws.Range("M5").Value = 1.048443520577824
$ws.Range("N5").Value = 1.036766464889325

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.030881910517546
$ws.Range("D6").Value = 1.038720617985422
$ws.Range("E6").Value = 0.9946255319796338
$ws.Range("F6").Value = 1.046079386021748
$ws.Range("I6").Value = 1.032520047006539
$ws.Range("J6").Value = 1.035324312456273
$ws.Range("K6").Value = 1.041137414832086
$ws.Range("L6").Value = 0.9971555583673453
$ws.Range("M6").Value = 1.0484782864519
$ws.Range("N6").Value = 1.03679459173332

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.030575452695149
$ws.Range("D7").Value = 1.038481919107974
$ws.Range("E7").Value = 0.9943035907982488
$ws.Range("F7").Value = 1.04578042703877
$ws.Range("I7").Value = 1.032474521963717
$ws.Range("J7").Value = 1.035131156476375
$ws.Range("K7").Value = 1.040958934199027
$ws.Range("L7").Value = 0.9968970624462087
$ws.Range("M7").Value = 1.048239223147946
$ws.Range("N7").Value = 1.036601161449775

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.029295096157606
$ws.Range("D8").Value = 1.037484316118142
$ws.Range("E8").Value = 0.9929600610674301
$ws.Range("F8").Value = 1.044531933346655
$ws.Range("I8").Value = 1.032281478896736
$ws.Range("J8").Value = 1.03432318442408
$ws.Range("K8").Value = 1.040211692379605
$ws.Range("L8").Value = 0.995817528259106
$ws.Range("M8").Value = 1.047239820307927
$ws.Range("N8").Value = 1.03579204198449

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.02704119834437
$ws.Range("D9").Value = 1.035726940582027
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.042336171814873
$ws.Range("I9").Value = 1.031930954979587
$ws.Range("J9").Value = 1.032897178665435
$ws.Range("K9").Value = 1.038890409372019
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.045478217784946
$ws.Range("N9").Value = 1.034364011134102

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.025540150082755
$ws.Range("D10").Value = 1.034555781838992
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.040875249968479
$ws.Range("I10").Value = 1.031690336014876
$ws.Range("J10").Value = 1.031945035105704
$ws.Range("K10").Value = 1.038006548382924
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.044303519805536
$ws.Range("N10").Value = 1.033410515421304

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.024890554008044
$ws.Range("D11").Value = 1.034048772305087
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.040243360491691
$ws.Range("I11").Value = 1.031584505739946
$ws.Range("J11").Value = 1.031532405849673
$ws.Range("K11").Value = 1.037623124490733
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.043794804158838
$ws.Range("N11").Value = 1.0329973001844

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.024649320807748
$ws.Range("D12").Value = 1.033860464018
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.040008754699573
$ws.Range("I12").Value = 1.031544949586817
$ws.Range("J12").Value = 1.031379085825585
$ws.Range("K12").Value = 1.037480598481516
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.043605835851357
$ws.Range("N12").Value = 1.032843762428292

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.024701063619933
$ws.Range("D13").Value = 1.033900855969242
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.040059073620852
$ws.Range("I13").Value = 1.031553445644755
$ws.Range("J13").Value = 1.031411975795112
$ws.Range("K13").Value = 1.037511175546644
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.043646370572475
$ws.Range("N13").Value = 1.032876699105348

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.024870612457424
$ws.Range("D14").Value = 1.034033206317291
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.040223965733483
$ws.Range("I14").Value = 1.031581241033476
$ws.Range("J14").Value = 1.031519733408212
$ws.Range("K14").Value = 1.037611345389875
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.043779184151168
$ws.Range("N14").Value = 1.032984609746618

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.024975084487258
$ws.Range("D15").Value = 1.034114754104023
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.04032557529553
$ws.Range("I15").Value = 1.031598334094703
$ws.Range("J15").Value = 1.031586119672427
$ws.Range("K15").Value = 1.037673049406546
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.043861013859027
$ws.Range("N15").Value = 1.033051090286944

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.025583269421472
$ws.Range("D16").Value = 1.034589432803082
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.040917201213748
$ws.Range("I16").Value = 1.031697325090032
$ws.Range("J16").Value = 1.031972412741841
$ws.Range("K16").Value = 1.038031980152488
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.044337280322963
$ws.Range("N16").Value = 1.033437931936825

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.025964866536373
$ws.Range("D17").Value = 1.034887216427733
$ws.Range("E17").Value = 0.989476357848556
$ws.Range("F17").Value = 1.041288500476899
$ws.Range("I17").Value = 1.031758980564254
$ws.Range("J17").Value = 1.032214632327234
$ws.Range("K17").Value = 1.038256939409886
$ws.Range("L17").Value = 0.9930127773699352
$ws.Range("M17").Value = 1.04463601328221
$ws.Range("N17").Value = 1.033680495501811

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.026187481126905
$ws.Range("D18").Value = 1.035060919110066
$ws.Range("E18").Value = 0.9897087662937556
$ws.Range("F18").Value = 1.041505140460765
$ws.Range("I18").Value = 1.031794784829934
$ws.Range("J18").Value = 1.032355881599927
$ws.Range("K18").Value = 1.038388086141917
$ws.Range("L18").Value = 0.9932001317071769
$ws.Range("M18").Value = 1.044810252959624
$ws.Range("N18").Value = 1.033821945364673

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.026263392948362
$ws.Range("D19").Value = 1.035120148988422
$ws.Range("E19").Value = 0.9897880325774034
$ws.Range("F19").Value = 1.041579020498628
$ws.Range("I19").Value = 1.03180696629122
$ws.Range("J19").Value = 1.032404038281155
$ws.Range("K19").Value = 1.038432792196618
$ws.Range("L19").Value = 0.9932640239640975
$ws.Range("M19").Value = 1.044869663107254
$ws.Range("N19").Value = 1.033870170433913

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.025923921107248
$ws.Range("D20").Value = 1.034855265965532
$ws.Range("E20").Value = 0.9894336180360679
$ws.Range("F20").Value = 1.041248656600344
$ws.Range("I20").Value = 1.031752381889093
$ws.Range("J20").Value = 1.03218864789903
$ws.Range("K20").Value = 1.038232810458636
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.044603962697975
$ws.Range("N20").Value = 1.03365447417274

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.024820683009246
$ws.Range("D21").Value = 1.03399423194298
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.04017540617915
$ws.Range("I21").Value = 1.031573062773864
$ws.Range("J21").Value = 1.031488002883539
$ws.Range("K21").Value = 1.037581850768166
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.043740074094524
$ws.Range("N21").Value = 1.032952834160961

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.02412735595172
$ws.Range("D22").Value = 1.033452968220656
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.039501225381426
$ws.Range("I22").Value = 1.031458894210074
$ws.Range("J22").Value = 1.031047183964253
$ws.Range("K22").Value = 1.037171957716408
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.043196863350135
$ws.Range("N22").Value = 1.032511389228247

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.024494871017845
$ws.Range("D23").Value = 1.033739892270093
$ws.Range("E23").Value = 0.9879432794643023
$ws.Range("F23").Value = 1.039858562802335
$ws.Range("I23").Value = 1.031519551919986
$ws.Range("J23").Value = 1.031280898212669
$ws.Range("K23").Value = 1.037389307223439
$ws.Range("L23").Value = 0.991776070289318
$ws.Range("M23").Value = 1.043484834053362
$ws.Range("N23").Value = 1.032745435377703

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.025942422474617
$ws.Range("D24").Value = 1.034869702971648
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.041266660124883
$ws.Range("I24").Value = 1.031755364035598
$ws.Range("J24").Value = 1.032200389246601
$ws.Range("K24").Value = 1.038243713503845
$ws.Range("L24").Value = 0.9929938892766442
$ws.Range("M24").Value = 1.044618444996602
$ws.Range("N24").Value = 1.033666232194371

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.027623614266691
$ws.Range("D25").Value = 1.03618119447623
$ws.Range("E25").Value = 0.9912096547607049
$ws.Range("F25").Value = 1.042903319097461
$ws.Range("I25").Value = 1.032022798397818
$ws.Range("J25").Value = 1.0332660981621
$ws.Range("K25").Value = 1.039232526255656
$ws.Range("L25").Value = 0.9944092447426414
$ws.Range("M25").Value = 1.045933690734839
$ws.Range("N25").Value = 1.034733454538768
